{"js": "// Troca o nome do ator \"Gestor\" por \"Ator\" em todas as ocorr\u00eancias do\n// corpo do documento (texto da descri\u00e7\u00e3o do caso de uso).\nconst body = context.document.body;\n\n// \"Gestor\" aparece como palavra isolada em 3 pontos do texto:\n//   - \"Gestor seleciona uma natureza financeira ...\"\n//   - \"Gestor altera as informa\u00e7\u00f5es desejadas.\"\n//   - \"Gestor clica no bot\u00e3o Salvar...\"\n// matchCase/matchWholeWord evitam tocar em qualquer outra ocorr\u00eancia\n// parcial (n\u00e3o existe nenhuma neste documento, mas \u00e9 mais seguro).\nconst results = body.search(\"Gestor\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (const result of results.items) {\n  result.insertText(\"Ator\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Troca o nome do ator \"Gestor\" por \"Ator\" em todas as ocorr\u00eancias do\n# corpo do documento (texto da descri\u00e7\u00e3o do caso de uso).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Gestor\"\n$find.Replacement.Text = \"Ator\"\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.MatchWildcards = $false\n\n# wdReplaceAll = 2 -> substitui todas as ocorr\u00eancias de uma s\u00f3 vez\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
